# Refresh the cryptocurrency price ("D") and 1h-volume-change ("E") columns
# with the latest scraped figures (GitHub Actions data refresh).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: target cell + new text. ForceText=$true locks the cell to Text
# format before the write and restores General afterwards, so a plain
# decimal-looking price string (e.g. "347.77") is kept as text instead of
# being auto-parsed into a numeric value by Excel.
$updates = @(
    @{ Cell = 'D2'; Text = '30.630.34'; ForceText = $false },
    @{ Cell = 'E2'; Text = '  +0.54%  '; ForceText = $false },
    @{ Cell = 'D3'; Text = '2.114.51'; ForceText = $false },
    @{ Cell = 'E3'; Text = '  +0.35%  '; ForceText = $false },
    @{ Cell = 'E4'; Text = '  +0.82%  '; ForceText = $false },
    @{ Cell = 'D5'; Text = '347.77'; ForceText = $true },
    @{ Cell = 'E5'; Text = '  +4.29%  '; ForceText = $false },
    @{ Cell = 'D6'; Text = '1.013'; ForceText = $true },
    @{ Cell = 'E6'; Text = '  +0.93%  '; ForceText = $false },
    @{ Cell = 'D7'; Text = '0.5271'; ForceText = $true },
    @{ Cell = 'E7'; Text = '  +0.55%  '; ForceText = $false },
    @{ Cell = 'D8'; Text = '0.4520'; ForceText = $true },
    @{ Cell = 'E8'; Text = '  -1.79%  '; ForceText = $false },
    @{ Cell = 'D9'; Text = '53.77'; ForceText = $true },
    @{ Cell = 'E9'; Text = '  +0.14%  '; ForceText = $false },
    @{ Cell = 'D10'; Text = '0.09048'; ForceText = $true },
    @{ Cell = 'E10'; Text = '  +1.08%  '; ForceText = $false },
    @{ Cell = 'D11'; Text = '1.174'; ForceText = $true },
    @{ Cell = 'E11'; Text = '  -0.28%  '; ForceText = $false },
    @{ Cell = 'D12'; Text = '24.47'; ForceText = $true },
    @{ Cell = 'E12'; Text = '  +0.29%  '; ForceText = $false },
    @{ Cell = 'D13'; Text = '2.107.22'; ForceText = $false },
    @{ Cell = 'E13'; Text = '  +0.19%  '; ForceText = $false },
    @{ Cell = 'D14'; Text = '6.826'; ForceText = $true },
    @{ Cell = 'E14'; Text = '  +0.76%  '; ForceText = $false },
    @{ Cell = 'D15'; Text = '8.086'; ForceText = $true },
    @{ Cell = 'E15'; Text = '  +2.99%  '; ForceText = $false },
    @{ Cell = 'D16'; Text = '99.97'; ForceText = $true },
    @{ Cell = 'E16'; Text = '  +3.59%  '; ForceText = $false },
    @{ Cell = 'E17'; Text = '  +4.52%  '; ForceText = $false },
    @{ Cell = 'D19'; Text = '0.06728'; ForceText = $true },
    @{ Cell = 'E19'; Text = '  +1.53%  '; ForceText = $false },
    @{ Cell = 'D20'; Text = '19.35'; ForceText = $true },
    @{ Cell = 'E20'; Text = '  +0.57%  '; ForceText = $false },
    @{ Cell = 'D21'; Text = '1.012'; ForceText = $true },
    @{ Cell = 'E21'; Text = '  +0.93%  '; ForceText = $false },
    @{ Cell = 'D22'; Text = '6.332'; ForceText = $true },
    @{ Cell = 'E22'; Text = '  +0.83%  '; ForceText = $false },
    @{ Cell = 'D23'; Text = '30.696.48'; ForceText = $false },
    @{ Cell = 'E23'; Text = '  +0.51%  '; ForceText = $false },
    @{ Cell = 'D24'; Text = '12.79'; ForceText = $true },
    @{ Cell = 'E24'; Text = '  +3.74%  '; ForceText = $false },
    @{ Cell = 'D25'; Text = '2.396'; ForceText = $true },
    @{ Cell = 'E25'; Text = '  +1.73%  '; ForceText = $false },
    @{ Cell = 'D26'; Text = '2.361.27'; ForceText = $false },
    @{ Cell = 'E26'; Text = '  +0.39%  '; ForceText = $false },
    @{ Cell = 'D27'; Text = '22.40'; ForceText = $true },
    @{ Cell = 'E27'; Text = '  +0.52%  '; ForceText = $false },
    @{ Cell = 'D28'; Text = '165.54'; ForceText = $true },
    @{ Cell = 'E28'; Text = '  +1.28%  '; ForceText = $false },
    @{ Cell = 'D29'; Text = '2.531'; ForceText = $true },
    @{ Cell = 'E29'; Text = '  -1.13%  '; ForceText = $false },
    @{ Cell = 'D30'; Text = '136.16'; ForceText = $true },
    @{ Cell = 'E30'; Text = '  +2.64%  '; ForceText = $false },
    @{ Cell = 'D31'; Text = '1.194'; ForceText = $true },
    @{ Cell = 'E31'; Text = '  +0.04%  '; ForceText = $false },
    @{ Cell = 'D32'; Text = '0.1074'; ForceText = $true },
    @{ Cell = 'E32'; Text = '  +0.24%  '; ForceText = $false },
    @{ Cell = 'D33'; Text = '1.636'; ForceText = $true },
    @{ Cell = 'E33'; Text = '  -2.85%  '; ForceText = $false },
    @{ Cell = 'D34'; Text = '6.360'; ForceText = $true },
    @{ Cell = 'E34'; Text = '  +3.51%  '; ForceText = $false },
    @{ Cell = 'D35'; Text = '4.008'; ForceText = $true },
    @{ Cell = 'E35'; Text = '  +1.93%  '; ForceText = $false },
    @{ Cell = 'D36'; Text = '5.927'; ForceText = $true },
    @{ Cell = 'E36'; Text = '  +6.88%  '; ForceText = $false },
    @{ Cell = 'D37'; Text = '10.25'; ForceText = $true },
    @{ Cell = 'E37'; Text = '  -1.78%  '; ForceText = $false },
    @{ Cell = 'D38'; Text = '0.02641'; ForceText = $true },
    @{ Cell = 'E38'; Text = '  +2.90%  '; ForceText = $false },
    @{ Cell = 'D39'; Text = '0.06844'; ForceText = $true },
    @{ Cell = 'E39'; Text = '  +0.41%  '; ForceText = $false },
    @{ Cell = 'D40'; Text = '0.2315'; ForceText = $true },
    @{ Cell = 'E40'; Text = '  +1.18%  '; ForceText = $false },
    @{ Cell = 'D41'; Text = '12.62'; ForceText = $true },
    @{ Cell = 'E41'; Text = '  -1.14%  '; ForceText = $false },
    @{ Cell = 'D42'; Text = '0.6883'; ForceText = $true },
    @{ Cell = 'E42'; Text = '  +0.14%  '; ForceText = $false },
    @{ Cell = 'E43'; Text = '  +1.88%  '; ForceText = $false },
    @{ Cell = 'D44'; Text = '14.83'; ForceText = $true },
    @{ Cell = 'E44'; Text = '  +6.57%  '; ForceText = $false },
    @{ Cell = 'D45'; Text = '2.326'; ForceText = $true },
    @{ Cell = 'E45'; Text = '  -0.91%  '; ForceText = $false },
    @{ Cell = 'D46'; Text = '0.6423'; ForceText = $true },
    @{ Cell = 'E46'; Text = '  +0.85%  '; ForceText = $false },
    @{ Cell = 'E47'; Text = '  +2.92%  '; ForceText = $false },
    @{ Cell = 'D48'; Text = '0.00000000360'; ForceText = $true },
    @{ Cell = 'E48'; Text = '  +2.86%  '; ForceText = $false },
    @{ Cell = 'D49'; Text = '1.256'; ForceText = $true },
    @{ Cell = 'E49'; Text = '  +0.90%  '; ForceText = $false },
    @{ Cell = 'D50'; Text = '82.90'; ForceText = $true },
    @{ Cell = 'E50'; Text = '  -0.18%  '; ForceText = $false },
    @{ Cell = 'D51'; Text = '0.07296'; ForceText = $true },
    @{ Cell = 'E51'; Text = '  +2.71%  '; ForceText = $false }
)

foreach ($u in $updates) {
    $range = $ws.Range($u.Cell)
    if ($u.ForceText) {
        $range.NumberFormat = "@"
        $range.Value = $u.Text
        $range.NumberFormat = "General"
    } else {
        $range.Value = $u.Text
    }
}
